$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the lab test search terms (rows 2-11) with the new lab + radiology
# test terms, per the "lab and radiology RL test" commit.
$ws.Range("A2").Value = "MRI Brain (CC, RADNET)"
$ws.Range("A3").Value = "CT Abdomen (CC, MIS)"
$ws.Range("A4").Value = "CT Chest (CC, MIS)"
$ws.Range("A5").Value = "US Kidneys - Bilat (CC, RADNET)"
$ws.Range("A6").Value = "US Thyroid (CC, MIS)"
$ws.Range("A7").Value = "NM Liver/Spleen Scan (CC, MIS)"
$ws.Range("A8").Value = "MM Mammography (CC, MIS)"
$ws.Range("A9").Value = "DX Shoulder Left (CC, MIS)"
$ws.Range("A10").Value = "DX Chest PA+Lateral (CC, MIS)"
$ws.Range("A11").Value = "Potassium Plasma Test, (CC, SOFTLAB, K5)"

# Move the active selection to A14 (was E26).
$ws.Range("A14").Select()
